# Auto-generated edit script: updates market-price derived columns (H:N)
# across the 8 Leve-profit worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 80
$ws.Range("H80").Value = 582.375
$ws.Range("I80").Value = 420
$ws.Range("J80").Value = 679.8
$ws.Range("K80").Value = 1260
$ws.Range("L80").Value = 2039.4
$ws.Range("M80").Value = -262
$ws.Range("N80").Value = -4035.4

# Row 82
$ws.Range("H82").Value = 16106.556
$ws.Range("I82").Value = 16106.556
$ws.Range("K82").Value = 48319.66800000001
$ws.Range("M82").Value = -47913.66800000001

# Row 83
$ws.Range("H83").Value = 582.375
$ws.Range("I83").Value = 420
$ws.Range("J83").Value = 679.8
$ws.Range("K83").Value = 3780
$ws.Range("L83").Value = 6118.2
$ws.Range("M83").Value = 1212
$ws.Range("N83").Value = -16102.2

# Row 85
$ws.Range("H85").Value = 16106.556
$ws.Range("I85").Value = 16106.556
$ws.Range("K85").Value = 48319.66800000001
$ws.Range("M85").Value = -46915.66800000001

# Row 98
$ws.Range("H98").Value = 5076.933
$ws.Range("I98").Value = 4781.154
$ws.Range("J98").Value = 6999.5
$ws.Range("K98").Value = 4781.154
$ws.Range("L98").Value = 6999.5
$ws.Range("M98").Value = -3283.154
$ws.Range("N98").Value = -9995.5

# Row 122
$ws.Range("H122").Value = 5076.933
$ws.Range("I122").Value = 4781.154
$ws.Range("J122").Value = 6999.5
$ws.Range("K122").Value = 14343.462
$ws.Range("L122").Value = 20998.5
$ws.Range("M122").Value = -11893.462
$ws.Range("N122").Value = -25898.5

# Row 131
$ws.Range("H131").Value = 2062.3333
$ws.Range("I131").Value = 2062.3333
$ws.Range("K131").Value = 6186.999899999999
$ws.Range("M131").Value = -1146.999899999999

# Row 137
$ws.Range("H137").Value = 1143.9231
$ws.Range("J137").Value = 1926.2142
$ws.Range("L137").Value = 5778.642599999999
$ws.Range("N137").Value = -10878.6426

# Row 141
$ws.Range("H141").Value = 65412.6
$ws.Range("I141").Value = 65412.6
$ws.Range("K141").Value = 196237.8
$ws.Range("M141").Value = -191057.8

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 1399438.2
$ws.Range("I32").Value = 1395270
$ws.Range("K32").Value = 1395270
$ws.Range("M32").Value = -1394983

# Row 45
$ws.Range("H45").Value = 2517.1667
$ws.Range("I45").Value = 992
$ws.Range("J45").Value = 4652.4
$ws.Range("K45").Value = 992
$ws.Range("L45").Value = 4652.4
$ws.Range("M45").Value = -615
$ws.Range("N45").Value = -5406.4

# Row 88
$ws.Range("H88").Value = 2452.2727
$ws.Range("I88").Value = 2016.8889
$ws.Range("J88").Value = 2753.6924
$ws.Range("K88").Value = 2016.8889
$ws.Range("L88").Value = 2753.6924
$ws.Range("M88").Value = -1610.8889
$ws.Range("N88").Value = -3565.6924

# Row 91
$ws.Range("H91").Value = 2452.2727
$ws.Range("I91").Value = 2016.8889
$ws.Range("J91").Value = 2753.6924
$ws.Range("K91").Value = 2016.8889
$ws.Range("L91").Value = 2753.6924
$ws.Range("M91").Value = -612.8888999999999
$ws.Range("N91").Value = -5561.6924

# Row 102
$ws.Range("H102").Value = 1201.3334
$ws.Range("I102").Value = 1241.6
$ws.Range("K102").Value = 1241.6
$ws.Range("M102").Value = 380.4000000000001

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 11191.583
$ws.Range("I94").Value = 3112.7896
$ws.Range("K94").Value = 3112.7896
$ws.Range("M94").Value = -2661.7896

# Row 134
$ws.Range("H134").Value = 20119618
$ws.Range("I134").Value = 9264095
$ws.Range("K134").Value = 27792285
$ws.Range("M134").Value = -27789750

$ws = $wb.Worksheets.Item("CRP")
# Row 56
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").Value = $null

# Row 74
$ws.Range("H74").Value = 26209
$ws.Range("J74").Value = 24313.5
$ws.Range("L74").Value = 24313.5
$ws.Range("N74").Value = -26061.5

# Row 77
$ws.Range("H77").Value = 26209
$ws.Range("J77").Value = 24313.5
$ws.Range("L77").Value = 72940.5
$ws.Range("N77").Value = -81676.5

# Row 103
$ws.Range("H103").Value = 25999.75
$ws.Range("J103").Value = 40000
$ws.Range("L103").Value = 40000
$ws.Range("N103").Value = -42344

$ws = $wb.Worksheets.Item("CUL")
# Row 81
$ws.Range("H81").Value = 3127.5715
$ws.Range("I81").Value = 2599
$ws.Range("K81").Value = 7797
$ws.Range("M81").Value = -6674

# Row 84
$ws.Range("H84").Value = 3127.5715
$ws.Range("I84").Value = 2599
$ws.Range("K84").Value = 23391
$ws.Range("M84").Value = -17775

# Row 108
$ws.Range("H108").Value = 3480
$ws.Range("I108").Value = 3480
$ws.Range("K108").Value = 10440
$ws.Range("M108").Value = -7560

# Row 122
$ws.Range("H122").Value = 1972.6666
$ws.Range("J122").Value = 1993.4286
$ws.Range("L122").Value = 17940.8574
$ws.Range("N122").Value = -22840.8574

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 4647.75
$ws.Range("J80").Value = 7264.4443
$ws.Range("L80").Value = 7264.4443
$ws.Range("N80").Value = -9260.444299999999

# Row 83
$ws.Range("H83").Value = 4647.75
$ws.Range("J83").Value = 7264.4443
$ws.Range("L83").Value = 36322.2215
$ws.Range("N83").Value = -46306.2215

# Row 107
$ws.Range("H107").Value = 592.875
$ws.Range("I107").Value = 314.42856
$ws.Range("K107").Value = 314.42856
$ws.Range("M107").Value = 1605.57144

# Row 113
$ws.Range("H113").Value = 3596.4546
$ws.Range("I113").Value = 3631.8
$ws.Range("K113").Value = 3631.8
$ws.Range("M113").Value = -1461.8

$ws = $wb.Worksheets.Item("LTW")
# Row 82
$ws.Range("H82").Value = 31510.611
$ws.Range("I82").Value = 4591.615
$ws.Range("K82").Value = 4591.615
$ws.Range("M82").Value = -4230.615

# Row 85
$ws.Range("H85").Value = 31510.611
$ws.Range("I85").Value = 4591.615
$ws.Range("K85").Value = 4591.615
$ws.Range("M85").Value = -3343.615

# Row 93
$ws.Range("H93").Value = 25761.277
$ws.Range("I93").Value = 2714
$ws.Range("J93").Value = 48808.555
$ws.Range("K93").Value = 2714
$ws.Range("L93").Value = 48808.555
$ws.Range("M93").Value = -1466
$ws.Range("N93").Value = -51304.555

$ws = $wb.Worksheets.Item("WVR")
# Row 20
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = $null
$ws.Range("N20").Value = $null

# Row 45
$ws.Range("H45").Value = 10700
$ws.Range("J45").Value = 10824
$ws.Range("L45").Value = 10824
$ws.Range("N45").Value = -11806

# Row 81
$ws.Range("H81").Value = 1661.25
$ws.Range("I81").Value = 1357.7273
$ws.Range("J81").Value = 5000
$ws.Range("K81").Value = 2715.4546
$ws.Range("L81").Value = 10000
$ws.Range("M81").Value = -1654.4546
$ws.Range("N81").Value = -12122

# Row 84
$ws.Range("H84").Value = 1661.25
$ws.Range("I84").Value = 1357.7273
$ws.Range("J84").Value = 5000
$ws.Range("K84").Value = 13577.273
$ws.Range("L84").Value = 50000
$ws.Range("M84").Value = -8273.273000000001
$ws.Range("N84").Value = -60608

# Row 100
$ws.Range("H100").Value = 447.63635
$ws.Range("I100").Value = 447.63635
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 895.2727
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -354.2727
$ws.Range("N100").Value = $null
